# Apply the scripted set of text replacements to the active document.
$d = $word.ActiveDocument

# Mapping of old text -> new text, in document order.
$replacements = @(
    @{ Old = "2025-02-17 Monday"; New = "2025-02-18 Tuesday" },
    @{ Old = "91×18="; New = "34×55=" },
    @{ Old = "43×12="; New = "63×57=" },
    @{ Old = "57×91="; New = "54×73=" },
    @{ Old = "48×14="; New = "45×95=" },
    @{ Old = "79×74="; New = "39×24=" },
    @{ Old = "12×64="; New = "33×94=" },
    @{ Old = "19×94="; New = "73×28=" },
    @{ Old = "93×93="; New = "47×20=" },
    @{ Old = "51×13="; New = "19×14=" },
    @{ Old = "30×42="; New = "70×56=" },
    @{ Old = "58×95="; New = "93×29=" },
    @{ Old = "58×21="; New = "97×92=" },
    @{ Old = "75×61="; New = "37×19=" },
    @{ Old = "69×94="; New = "24×85=" },
    @{ Old = "94×78="; New = "84×82=" },
    @{ Old = "31×23="; New = "16×77=" },
    @{ Old = "85×45="; New = "41×41=" },
    @{ Old = "36×21="; New = "83×91=" },
    @{ Old = "34×35="; New = "59×58=" },
    @{ Old = "33×95="; New = "47×72=" },
    @{ Old = "20×45="; New = "44×33=" },
    @{ Old = "75×33="; New = "29×14=" },
    @{ Old = "82×82="; New = "72×93=" },
    @{ Old = "16×18="; New = "94×75=" },
    @{ Old = "27×27="; New = "42×33=" }
)

foreach ($item in $replacements) {
    $range = $d.Content
    $range.Find.Execute($item.Old, $true, $true, $false, $false, $false, $true, 1, $false, $item.New, 2)
}

$d.Save()
